# ---------------------------------------------------------------------------
# Edit: Thu, Jun 11, 2020  7:05:06 AM
#
# 1) The table on slide 5 gets a new built-in table style GUID applied.
# 2) The presentation's active theme colour palette is switched from the
#    "Integral / Red Violet" palette to the default "Office Theme / Office"
#    palette (i.e. the Design gallery selection was changed back to the
#    stock Office theme).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style -> Medium Style 2 - Accent 1 -----------------------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{926FA7BD-E4F3-4246-AC09-6E23D763A1EB}")
    }
}

# --- 2) Swap the applied colour theme back to the stock "Office" palette ----
function ConvertTo-ComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (in clrScheme document order)
$officeThemeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47",
    "0563C1", "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-ComRgb($officeThemeColors[$i - 1])
}
